$d = $word.ActiveDocument

# --- Interfaces block: re-derive the type hierarchy chain ---

# 1. Value : ID;  ->  Context: ID;
$d.Content.Find.Execute("Value : ID;", $false, $false, $false, $false, $false, $true, 1, $false, "Context: ID;", 2)

# 2. Sign : Value;  ->  Object : Context;
$d.Content.Find.Execute("Sign : Value;", $false, $false, $false, $false, $false, $true, 1, $false, "Object : Context;", 2)

# 3. Object : Sign;  ->  Sign : Object;
$d.Content.Find.Execute("Object : Sign;", $false, $false, $false, $false, $false, $true, 1, $false, "Sign : Object;", 2)

# 4. Context: Object;  ->  Value : Sign;
$d.Content.Find.Execute("Context: Object;", $false, $false, $false, $false, $false, $true, 1, $false, "Value : Sign;", 2)

# 5. Context: (Context, Object, Sign, Value);  ->  Context (Context, Object, Sign, Value);
$d.Content.Find.Execute("Context: (Context, Object, Sign, Value);", $false, $false, $false, $false, $false, $true, 1, $false, "Context (Context, Object, Sign, Value);", 2)

# 6. Resource: (...)  ->  Resource : Value (...)
$d.Content.Find.Execute("Resource: (Resource, Resource, Resource, Resource);", $false, $false, $false, $false, $false, $true, 1, $false, "Resource : Value (Resource, Resource, Resource, Resource);", 2)

# 7. Statement: (...)  ->  Statement : Resource (...)
$d.Content.Find.Execute("Statement: (Statement, Resource, Resource, Resource);", $false, $false, $false, $false, $false, $true, 1, $false, "Statement : Resource (Statement, Resource, Resource, Resource);", 2)

# 8. Role: (...)  ->  Role : Statement (...)
$d.Content.Find.Execute("Role: (Statement, Role, Resource, Resource);", $false, $false, $false, $false, $false, $true, 1, $false, "Role : Statement (Statement, Role, Resource, Resource);", 2)

# --- Relationship line: mark it with a footnote-style "*" and add the explanatory note ---

# 9. Append " *" to the Relationship line
$d.Content.Find.Execute("Relationship : Relation (Relationship, Role, Kind, Relation);", $false, $false, $false, $false, $false, $true, 1, $false, "Relationship : Relation (Relationship, Role, Kind, Relation); *", 2)

$target = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Relationship : Relation (Relationship, Role, Kind, Relation); *`r") {
        $target = $i
        break
    }
}

$p = $d.Paragraphs.Item($target)

# Insert a blank paragraph, then the footnote-text paragraph, right after it
$p.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Item($target + 1)
$pBlank.Range.InsertParagraphAfter()
$pNote = $d.Paragraphs.Item($target + 2)
$pNote.Range.Text = "*: Relationship: Aggregated Relation Statement Relation (Object) Roles / Kinds."

# Only the original "Relationship" paragraph gets the new (0,0) explicit indent
$p.LeftIndent = 0
$p.FirstLineIndent = 0

# --- Misc text updates ---

# 10. Context Monad sentence gets appended text
$d.Content.Find.Execute("Context Monad & type hierarchy (AST).", $false, $false, $false, $false, $false, $true, 1, $false, "Context Monad & type hierarchy (AST). CSPO parameterized types & aggregation (layers hierarchies specializations).", 2)
